# Apply cryptos list update (commit: "Updated cryptos list on Fri Jul 26 19:14:45 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '67.670.00'
$ws.Range("E2").Value = "'" + '  +4.68%  '
$ws.Range("D3").Value = "'" + '3.268.83'
$ws.Range("E3").Value = "'" + '  +4.98%  '
$ws.Range("E4").Value = "'" + '  +0.01%  '
$ws.Range("D5").Value = "'" + '580.25'
$ws.Range("E5").Value = "'" + '  +2.79%  '
$ws.Range("D6").Value = "'" + '182.80'
$ws.Range("E6").Value = "'" + '  +9.48%  '
$ws.Range("E7").Value = "'" + '  +0.12%  '
$ws.Range("D8").Value = "'" + '0.601'
$ws.Range("E8").Value = "'" + '  +0.25%  '
$ws.Range("D9").Value = "'" + '3.268.56'
$ws.Range("E9").Value = "'" + '  +5.09%  '
$ws.Range("D10").Value = "'" + '0.133'
$ws.Range("E10").Value = "'" + '  +9.72%  '
$ws.Range("D11").Value = "'" + '6.74'
$ws.Range("E11").Value = "'" + '  +3.98%  '
$ws.Range("D12").Value = "'" + '0.418'
$ws.Range("E12").Value = "'" + '  +8.60%  '
$ws.Range("D13").Value = "'" + '3.831.69'
$ws.Range("E13").Value = "'" + '  +4.98%  '
$ws.Range("D15").Value = "'" + '28.61'
$ws.Range("E15").Value = "'" + '  +8.48%  '
$ws.Range("D16").Value = "'" + '67.631.71'
$ws.Range("E16").Value = "'" + '  +4.78%  '
$ws.Range("D17").Value = "'" + '0.0000169'
$ws.Range("E17").Value = "'" + '  +5.99%  '
$ws.Range("D18").Value = "'" + '3.260.16'
$ws.Range("E18").Value = "'" + '  +4.74%  '
$ws.Range("D19").Value = "'" + '5.86'
$ws.Range("E19").Value = "'" + '  +4.40%  '
$ws.Range("D20").Value = "'" + '13.62'
$ws.Range("E20").Value = "'" + '  +8.51%  '
$ws.Range("D21").Value = "'" + '375.48'
$ws.Range("E21").Value = "'" + '  +6.99%  '
$ws.Range("D22").Value = "'" + '7.65'
$ws.Range("E22").Value = "'" + '  +7.28%  '
$ws.Range("E23").Value = "'" + '  +0.14%  '
$ws.Range("D24").Value = "'" + '71.40'
$ws.Range("E24").Value = "'" + '  +4.34%  '
$ws.Range("D25").Value = "'" + '0.514'
$ws.Range("E25").Value = "'" + '  +5.49%  '
$ws.Range("D26").Value = "'" + '0.0000120'
$ws.Range("E26").Value = "'" + '  +6.59%  '
$ws.Range("D27").Value = "'" + '9.65'
$ws.Range("E27").Value = "'" + '  +1.11%  '
$ws.Range("D28").Value = "'" + '0.181'
$ws.Range("E28").Value = "'" + '  +3.68%  '
$ws.Range("E29").Value = "'" + '  -0.19%  '
$ws.Range("D30").Value = "'" + '5.74'
$ws.Range("E30").Value = "'" + '  +10.56%  '
$ws.Range("E31").Value = "'" + '  +4.74%  '
$ws.Range("D32").Value = "'" + '22.76'
$ws.Range("E32").Value = "'" + '  +5.81%  '
$ws.Range("E33").Value = "'" + '  +0.05%  '
$ws.Range("D34").Value = "'" + '1.28'
$ws.Range("E34").Value = "'" + '  +9.40%  '
$ws.Range("D35").Value = "'" + '6.95'
$ws.Range("E35").Value = "'" + '  +7.21%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'" + '1.51'
$ws.Range("E36").Value = "'" + '  +7.18%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = "'" + '163.35'
$ws.Range("E37").Value = "'" + '  +3.51%  '
$ws.Range("D38").Value = "'" + '0.854'
$ws.Range("E38").Value = "'" + '  +4.02%  '
$ws.Range("E39").Value = "'" + '  +6.58%  '
$ws.Range("D40").Value = "'" + '6.86'
$ws.Range("E40").Value = "'" + '  +13.50%  '
$ws.Range("D41").Value = "'" + '4.70'
$ws.Range("E41").Value = "'" + '  +14.75%  '
$ws.Range("D42").Value = "'" + '26.90'
$ws.Range("E42").Value = "'" + '  +4.99%  '
$ws.Range("D43").Value = "'" + '2.63'
$ws.Range("E43").Value = "'" + '  +10.65%  '
$ws.Range("D44").Value = "'" + '358.82'
$ws.Range("E44").Value = "'" + '  +14.30%  '
$ws.Range("D45").Value = "'" + '2.710.36'
$ws.Range("E45").Value = "'" + '  +3.25%  '
$ws.Range("D46").Value = "'" + '25.53'
$ws.Range("E46").Value = "'" + '  +8.63%  '
$ws.Range("D47").Value = "'" + '40.92'
$ws.Range("E47").Value = "'" + '  +4.48%  '
$ws.Range("D48").Value = "'" + '0.0684'
$ws.Range("E48").Value = "'" + '  +6.57%  '
$ws.Range("D49").Value = "'" + '0.0281'
$ws.Range("E49").Value = "'" + '  +5.28%  '
$ws.Range("E50").Value = "'" + '  +8.98%  '
$ws.Range("D51").Value = "'" + '0.103'
$ws.Range("E51").Value = "'" + '  +1.51%  '
